$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column R into the new column S (row by row) so that
# each new S cell inherits the same style as its R counterpart.
$ws.Range("R1:R20").Copy()
$ws.Range("S1:S20").PasteSpecial(-4122)

# New header for column S (mirrors the T16 header pattern in R1)
$ws.Range("S1").Value = "T17: 5/4/2020"

# New data values for column S, rows 2-19 (per-department counts for T17)
$ws.Range("S2").Value = 11
$ws.Range("S3").Value = 2
$ws.Range("S4").Value = 20
$ws.Range("S5").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("S7").Value = 195
$ws.Range("S8").Value = 0
$ws.Range("S9").Value = 47
$ws.Range("S10").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("S13").Value = 1
$ws.Range("S14").Value = 4
$ws.Range("S15").Value = 0
$ws.Range("S16").Value = 0
$ws.Range("S17").Value = 8
$ws.Range("S18").Value = 0
$ws.Range("S19").Value = 8

# Sum row
$ws.Range("S20").Formula = "=SUM(S2:S19)"

# Column width for the new column (matches the target width for col S)
$ws.Columns.Item(19).ColumnWidth = 14.45

# Update view: select S5 as the active cell (matches the new selection)
$ws.Range("S5").Select() | Out-Null
